# Applies schema changes to UT IAC London daily hearing list:
# "Applicant" column header becomes "Case title", and the
# per-row applicant names become case names ("Case A"/"Case B").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Case title"
$ws.Range("B2").Value = "Case A"
$ws.Range("B3").Value = "Case B"

# Update the active selection to match the saved state of the workbook.
$ws.Range("G11").Select()
